$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new shared strings / rows -------------------------------------
# Row 23: new feature "Tab Bar" (highlighted with the new light-green fill)
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "Tab Bar"
$ws.Range("C23").Value = 10
$ws.Range("D23").Value = 4

# Row 24: new feature "Checkout" (no special fill)
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "Checkout"
$ws.Range("C24").Value = 10
$ws.Range("D24").Value = 8

# --- Apply the existing yellow highlight to rows 2-8 --------------------
$ws.Range("A2:D8").Interior.Color = 65535

# --- Apply a new light-green highlight to row 12 and row 23 -------------
$ws.Range("A12:D12").Interior.Color = 14348258
$ws.Range("A23:D23").Interior.Color = 14348258

# --- Update selection / active cell to match the new state --------------
$ws.Range("A24:D24").Select() | Out-Null
